# lesson #27 video record uploaded into youtube
# Fill in the row for lesson #27 (row 30): lesson name, hours, date, YouTube link.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (Lesson name) - copy formatting from the row above (C29) then set the text.
$ws.Range("C29").Copy()
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C30").Value = "Step Project #1, Algorithms #3"

# Column D (amount of hours)
$ws.Range("D30").Value = 2

# Column E (date) - copy formatting from the row above (E29) then set the date value.
$ws.Range("E29").Copy()
$ws.Range("E30").PasteSpecial(-4122)
$ws.Range("E30").Value = Get-Date -Year 2021 -Month 1 -Day 6 -Hour 0 -Minute 0 -Second 0

# Column F (YouTube link)
$ws.Range("F30").Value = "https://youtu.be/uJR7GhziAKs "

# Row height grew slightly to fit the new content.
$ws.Rows(30).RowHeight = 14.3

# Move the active selection to the newly filled cell.
$ws.Range("E30").Select()
